# Insert a new weekly Apio (celery) price record for Femacal de La Calera,
# dated 2021-09-09 (serial 44448), ahead of the existing 2020-12-02 record,
# shifting the subsequent rows (old 222-226) down to 223-227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 222; existing rows 222:226 shift down to 223:227,
# carrying their formatting (incl. the date style on column D) with them.
$ws.Rows.Item(222).Insert()

$ws.Range("A222").Value = 3
$ws.Range("B222").Value = "Femacal de La Calera"
$ws.Range("C222").Value = "Coquimbo"
$ws.Range("D222").Value = 44448
$ws.Range("E222").Value = 5
$ws.Range("F222").Value = 100112017
$ws.Range("G222").Value = "Apio"
$ws.Range("H222").Value = "Americana (o)"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 160
$ws.Range("K222").Value = 9000
$ws.Range("L222").Value = 9000
$ws.Range("M222").Value = 9000
$ws.Range("N222").Value = "$/docena de matas"
$ws.Range("O222").Value = "Pan de Azúcar"
$ws.Range("P222").Value = 1500
$ws.Range("Q222").Value = 6
$ws.Range("R222").Value = "Hortaliza"
